$d = $word.ActiveDocument

# The "Requisitos" section ends with the paragraph "LOM3016: Introdução à
# Ciência dos Materiais (Requisito)". Immediately after it, the document
# has an empty paragraph, then a "Ver no Jupiter Salvar em pdf Salvar em
# docx" paragraph, then a "© 2020 ... Creative Commons Attribution"
# paragraph. All three of those paragraphs (site-footer boilerplate) are
# removed, while the requirement paragraph itself and the remaining final
# empty paragraph (right before the page-break paragraph) are kept intact.

$reqPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOM3016: Introdução à Ciência dos Materiais (Requisito)*") {
        $reqPara = $p
    }
}

$copyrightPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Powered by Jekyll and Github pages*") {
        $copyrightPara = $p
    }
}

$delStart = $reqPara.Range.End
$delEnd = $copyrightPara.Range.End

$rangeToDelete = $d.Range($delStart, $delEnd)
$rangeToDelete.Delete()
